$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-10 with new values
$ws.Range("A2").Value = "Q0"
$ws.Range("B2").Value = 0.1500504596286886
$ws.Range("C2").Value = 0.9002586435756988
$ws.Range("D2").Value = 3.690897970745633
$ws.Range("E2").Value = 1.921170989460759
$ws.Range("F2").Value = 1.934360485255287
$ws.Range("G2").Value = 51

$ws.Range("A3").Value = "Q1"
$ws.Range("B3").Value = 0.1741937731364769
$ws.Range("C3").Value = 0.861216386388353
$ws.Range("D3").Value = 3.555569236642697
$ws.Range("E3").Value = 1.885621710906696
$ws.Range("F3").Value = 1.896620454657085
$ws.Range("G3").Value = 50

$ws.Range("A4").Value = "Q2"
$ws.Range("B4").Value = 0.1975424125742846
$ws.Range("C4").Value = 0.9344996350649757
$ws.Range("D4").Value = 3.868064466556859
$ws.Range("E4").Value = 1.966739552293811
$ws.Range("F4").Value = 1.977071864882974
$ws.Range("G4").Value = 49

$ws.Range("A5").Value = "Q3"
$ws.Range("B5").Value = 0.1696732928241277
$ws.Range("C5").Value = 0.8834793121284986
$ws.Range("D5").Value = 3.720086392232908
$ws.Range("E5").Value = 1.928752548211643
$ws.Range("F5").Value = 1.941606450278617
$ws.Range("G5").Value = 48

$ws.Range("A6").Value = "Q4"
$ws.Range("B6").Value = 0.1955713809712057
$ws.Range("C6").Value = 0.8996788932259349
$ws.Range("D6").Value = 3.907687536052267
$ws.Range("E6").Value = 1.976787175204318
$ws.Range("F6").Value = 1.98835550598802
$ws.Range("G6").Value = 47

$ws.Range("A7").Value = "Q5"
$ws.Range("B7").Value = 0.1353573482079271
$ws.Range("C7").Value = 0.8312788515306436
$ws.Range("D7").Value = 3.56529404348522
$ws.Range("E7").Value = 1.888198623949615
$ws.Range("F7").Value = 1.904151790526757
$ws.Range("G7").Value = 46

$ws.Range("A8").Value = "Q6"
$ws.Range("B8").Value = 0.1150417487834175
$ws.Range("C8").Value = 0.8182950210794
$ws.Range("D8").Value = 3.666196212993372
$ws.Range("E8").Value = 1.914731368363033
$ws.Range("F8").Value = 1.932869230905421
$ws.Range("G8").Value = 45

$ws.Range("A9").Value = "Q7"
$ws.Range("B9").Value = 0.1478860465117961
$ws.Range("C9").Value = 0.8234762687095699
$ws.Range("D9").Value = 3.667743996575232
$ws.Range("E9").Value = 1.915135503450143
$ws.Range("F9").Value = 1.931492033276065
$ws.Range("G9").Value = 44

$ws.Range("A10").Value = "Q8"
$ws.Range("B10").Value = 0.1233280604110242
$ws.Range("C10").Value = 0.8107235582056446
$ws.Range("D10").Value = 3.724587762219358
$ws.Range("E10").Value = 1.929919107688029
$ws.Range("F10").Value = 1.948767937542836
$ws.Range("G10").Value = 43

# New row 11
$ws.Range("A11").Value = "Q9"
$ws.Range("B11").Value = 0.0913265010797303
$ws.Range("C11").Value = 0.785852673267451
$ws.Range("D11").Value = 3.688102942804829
$ws.Range("E11").Value = 1.92044342348449
$ws.Range("F11").Value = 1.941523297764307
$ws.Range("G11").Value = 42

# Copy style from A10 (header-like style) to A11
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122) # xlPasteFormats
